# PIB_Nacional.xlsx edit script
# - Replace quarterly date headers in column B with "Trim-<roman>-<year>" text labels
# - Drop the oldest 4 rows (2017 data)
# - Apply a one-decimal thousands number format to the GDP figures
# - Right-align the new period labels
# - Fix a typo in the units subtitle
# - Widen column B slightly to fit the new labels
# - Rename the worksheet tab

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: "aprecios" -> "a precios"
$ws.Range("B3").Value = "(Millones de pesos anulizados  a precios de 2018)"

# Quarter labels for rows 5..29 (newest first, matches existing row order)
$labels = @(
  "Trim-I-2024","Trim-IV-2023","Trim-III-2023","Trim-II-2023","Trim-I-2023",
  "Trim-IV-2022","Trim-III-2022","Trim-II-2022","Trim-I-2022",
  "Trim-IV-2021","Trim-III-2021","Trim-II-2021","Trim-I-2021",
  "Trim-IV-2020","Trim-III-2020","Trim-II-2020","Trim-I-2020",
  "Trim-IV-2019","Trim-III-2019","Trim-II-2019","Trim-I-2019",
  "Trim-IV-2018","Trim-III-2018","Trim-II-2018","Trim-I-2018"
)

for ($i = 0; $i -lt $labels.Count; $i++) {
  $row = 5 + $i
  $ws.Cells.Item($row, 2).Value = $labels[$i]
}

# Drop the four oldest rows (2017 Q4..Q1) that no longer appear in the table
$ws.Rows("30:33").Delete()

# Right align the new text period labels (column B, data rows)
$ws.Range("B5:B29").HorizontalAlignment = -4152

# One-decimal thousands format for the GDP values
$ws.Range("C5:C29").NumberFormat = "#,##0.0"

# Column B is a bit wider now that it holds text labels instead of dates
$ws.Columns("B").ColumnWidth = 11.19921875

# Rename the sheet tab
$ws.Name = "C_1.1"

# Reset the active selection to the top-left cell
$ws.Range("B1").Select()
